$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.713.43'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.601.63'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.55'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.11'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.22%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0603'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = '1.830.72'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').Value = '1.604.78'
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.552'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '29.724.24'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.03'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.82'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.91'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.72%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.58'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.45'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0478'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.19'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.18%  '
$ws.Range('D34').Value = '1.419.36'
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.89'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '55.82'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0494'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.815'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.30'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.989'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +18.23%  '
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Value = '1.740.01'
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0524'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.60%  '
